# Fruta / hortaliza, semanal
# Insert the latest weekly price record for
# "Agrícola del Norte S.A. de Arica - Plátano" at the top of this
# sub-block (row 133), pushing all subsequent rows down by one.
# The new row duplicates the (former) first record's data but carries
# the newest survey date (serial 44543 = 2021-12-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 133 (and everything below it) down by one row.
$ws.Rows.Item(133).Insert()

# Copy the data that is now on row 134 (the original row 133 content)
# into the freshly inserted, still-empty row 133.
$srcValues = $ws.Range("A134:T134").Value()
$ws.Range("A133:T133").Value = $srcValues

# The new row records a later survey date than the row it was copied
# from; every other field is unchanged.
$ws.Range("D133").Value = 44543
